# Update vm_pu results for the case with 380 kV (Case_2_255).
# Rows 2-25 correspond to bus indices 0-23; columns B..F and I..N hold
# the per-bus voltage magnitudes (p.u.). The slack/reference voltage (B)
# drops from 1.05 to 1.02 p.u., and the downstream bus voltages shift
# accordingly to their newly recomputed load-flow values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$columns = @("B", "C", "D", "E", "F", "I", "J", "K", "L", "M", "N")

$newValues = @{
    2 = @(1.02, 1.040678082160796, 1.048405537223276, 1.054924486636589, 1.06109912910135, 1.043560072786564, 1.045763392179102, 1.051165473218208, 1.057666357619507, 1.063824110846303, 1.047248496146734)
    3 = @(1.02, 1.041518597462012, 1.049074876947421, 1.055820316700731, 1.061969849354246, 1.043781141583781, 1.046250167615262, 1.051647221153132, 1.058375316630267, 1.064509244167361, 1.047735962859848)
    4 = @(1.02, 1.042062869592612, 1.049508332350636, 1.056401285435298, 1.062534254114536, 1.043923036678041, 1.046564874501037, 1.051958609090521, 1.058834729077725, 1.064952917871863, 1.048051116665497)
    5 = @(1.02, 1.042291776450875, 1.049690638600843, 1.056645835768056, 1.062771765210076, 1.043982413009736, 1.046697111606839, 1.052089435093248, 1.059028024633009, 1.065139519761261, 1.048183541563158)
    6 = @(1.02, 1.042330216438047, 1.049721253358664, 1.056686915035278, 1.062811658115873, 1.04399236632861, 1.046719310913582, 1.0521113965403, 1.059060489062086, 1.065170855772477, 1.048205772395463)
    7 = @(1.02, 1.042065927886431, 1.049510768017155, 1.056404551908713, 1.062537426826976, 1.043923831154377, 1.046566641719232, 1.051960357515715, 1.058837311281225, 1.064955410937251, 1.048052886393345)
    8 = @(1.02, 1.040962053987035, 1.048631670686389, 1.055226964959448, 1.061393186907888, 1.043635021924158, 1.045927955281886, 1.051328350845638, 1.057905814309049, 1.064055582173175, 1.047413292948)
    9 = @(1.02, 1.039020036834261, 1.047085326957207, 1.053161992944773, 1.059384555928722, 1.043117320229542, 1.044800497291591, 1.050212162016438, 1.056269584019701, 1.062472690935532, 1.046284233838037)
    10 = @(1.02, 1.037727565982495, 1.046056371737986, 1.051792232048076, 1.058050729991408, 1.042766326131914, 1.044047576731498, 1.049466420707327, 1.055182339745087, 1.061419350700095, 1.045530244044388)
    11 = @(1.02, 1.037168454086313, 1.045611304038756, 1.051200766199533, 1.057474439134542, 1.04261296151086, 1.04372126447029, 1.049143138810554, 1.054712417841726, 1.060963717311364, 1.04520346838234)
    12 = @(1.02, 1.036960857124901, 1.045446059386706, 1.050981319086591, 1.057260570858467, 1.042555788196706, 1.043600014763423, 1.049023003016136, 1.054537998941925, 1.060794546904589, 1.045082046486982)
    13 = @(1.02, 1.037005383639133, 1.045481501581062, 1.051028379911753, 1.057306437620772, 1.04256806141581, 1.043626025155332, 1.049048774998466, 1.054575406434961, 1.060830831242976, 1.045108093816631)
    14 = @(1.02, 1.037151292373479, 1.045597643359046, 1.051182621529789, 1.057456756792788, 1.04260823976256, 1.043711242798121, 1.049133209453042, 1.054697997639608, 1.060949732152696, 1.045193432478248)
    15 = @(1.02, 1.03724120248271, 1.045669211906148, 1.051277688014491, 1.057549398866609, 1.042632967578652, 1.04376374256342, 1.049185225121791, 1.05477354753309, 1.061023000529774, 1.045246006799235)
    16 = @(1.02, 1.037764683828159, 1.046085919603158, 1.051831520598165, 1.058089003316053, 1.042776475385942, 1.044069226934756, 1.049487868151056, 1.055213545204956, 1.061449599617684, 1.045551924993417)
    17 = @(1.02, 1.038093194533921, 1.046347438151651, 1.052179368073657, 1.058427822778089, 1.042866124603185, 1.044260771674087, 1.049677609673823, 1.05548977597619, 1.061717320833421, 1.045743741748248)
    18 = @(1.02, 1.038284860936406, 1.046500023377989, 1.052382420668677, 1.058625572329136, 1.042918282095403, 1.04437246816107, 1.049788246816154, 1.055650979785111, 1.061873523323712, 1.045855596857055)
    19 = @(1.02, 1.038350222944773, 1.046552058727898, 1.052451683278859, 1.058693020435257, 1.042936043806309, 1.044410548937703, 1.04982596509957, 1.05570596016297, 1.061926791962426, 1.045893731712759)
    20 = @(1.02, 1.038057943089413, 1.046319374928444, 1.052142030867524, 1.058391458035455, 1.042856519881299, 1.044240223657016, 1.049657255906221, 1.055460130425102, 1.061688592190811, 1.045723164550637)
    21 = @(1.02, 1.037108323609821, 1.045563440462969, 1.051137194311986, 1.057412486233797, 1.042596413945378, 1.043686149498521, 1.049108347105312, 1.054661893956495, 1.060914716795162, 1.045168303543284)
    22 = @(1.02, 1.036511735845398, 1.045088579434233, 1.050506858906994, 1.056798078446814, 1.042431678530056, 1.043337534012247, 1.048762911624089, 1.054160769784984, 1.060428568643917, 1.044819192983032)
    23 = @(1.02, 1.036827952719895, 1.045340271404224, 1.050840874041588, 1.057123681756243, 1.042519121055333, 1.043522364706541, 1.048946062929764, 1.054426352763277, 1.060686244843437, 1.045004286158113)
    24 = @(1.02, 1.038073871540329, 1.046332055349484, 1.052158901452603, 1.058407889320945, 1.042860860253889, 1.044249508509085, 1.049666453008681, 1.055473525717058, 1.061701573288209, 1.04573246258826)
    25 = @(1.02, 1.039521712341908, 1.047484759122531, 1.053694631529228, 1.059902915681178, 1.043252194664611, 1.045092203283938, 1.050501013710363, 1.056691965000601, 1.062881574370362, 1.046576354086352)
}

foreach ($row in $newValues.Keys) {
    $rowValues = $newValues[$row]
    for ($i = 0; $i -lt $columns.Count; $i++) {
        $ws.Range($columns[$i] + $row).Value = $rowValues[$i]
    }
}
